$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.65"
$ws.Range("E2").Value = "'0.44%"
$ws.Range("D3").Value = "'36.13"
$ws.Range("E3").Value = "'-2.46%"
$ws.Range("D4").Value = "'5.058"
$ws.Range("E4").Value = "'0.52%"
$ws.Range("D5").Value = "'0.07863"
$ws.Range("E5").Value = "'0.41%"
$ws.Range("D6").Value = "'2.140"
$ws.Range("E6").Value = "'-3.16%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.138"
$ws.Range("E7").Value = "'2.80%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.936"
$ws.Range("E8").Value = "'-0.70%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9220"
$ws.Range("E9").Value = "'-0.09%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09744"
$ws.Range("E10").Value = "'-1.38%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1849"
$ws.Range("E11").Value = "'-1.74%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08625"
$ws.Range("E12").Value = "'0.32%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03588"
$ws.Range("E13").Value = "'-0.84%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09920"
$ws.Range("E14").Value = "'-0.18%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001432"
$ws.Range("E15").Value = "'-3.87%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005686"
$ws.Range("E16").Value = "'-0.71%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.471"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("D18").Value = "'2.752"
$ws.Range("E18").Value = "'17.31%"
$ws.Range("D19").Value = "'0.3376"
$ws.Range("E19").Value = "'-1.67%"
$ws.Range("D20").Value = "'0.1349"
$ws.Range("E20").Value = "'1.67%"
$ws.Range("D21").Value = "'5.149"
$ws.Range("E21").Value = "'7.87%"
$ws.Range("D22").Value = "'0.2212"
$ws.Range("E22").Value = "'0.49%"
$ws.Range("D23").Value = "'0.04564"
$ws.Range("E23").Value = "'-0.99%"
$ws.Range("D24").Value = "'0.001233"
$ws.Range("E24").Value = "'-1.46%"
$ws.Range("D25").Value = "'0.004798"
$ws.Range("E25").Value = "'-7.65%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("D27").Value = "'0.0004760"
$ws.Range("E27").Value = "'75.11%"
$ws.Range("D39").Value = "'0.01856"
$ws.Range("E39").Value = "'0.75%"
$ws.Range("D40").Value = "'0.04721"
$ws.Range("E40").Value = "'-0.72%"
$ws.Range("D41").Value = "'0.007790"
$ws.Range("E41").Value = "'-2.07%"
$ws.Range("D42").Value = "'0.1382"
$ws.Range("E42").Value = "'-1.53%"
$ws.Range("D43").Value = "'0.007770"
$ws.Range("E43").Value = "'3.11%"
$ws.Range("D44").Value = "'0.002162"
$ws.Range("E44").Value = "'-3.60%"
$ws.Range("D45").Value = "'0.01134"
$ws.Range("E45").Value = "'8.91%"
$ws.Range("D46").Value = "'0.00006385"
$ws.Range("E46").Value = "'1.46%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("E48").Value = "'0.15%"
$ws.Range("D49").Value = "'53.11"
$ws.Range("E49").Value = "'48.38%"
$ws.Range("D50").Value = "'0.001905"
$ws.Range("E50").Value = "'-29.21%"
$ws.Range("D51").Value = "'0.00002105"
$ws.Range("E51").Value = "'0.13%"
